$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 3204.4546
$ws.Range("I138").Value = 2935.1
$ws.Range("J138").Value = 5898
$ws.Range("K138").Value = 8805.299999999999
$ws.Range("L138").Value = 17694
$ws.Range("M138").Value = -3665.299999999999
$ws.Range("N138").Value = -27974

$ws = $wb.Worksheets.Item(2)
$ws.Range("N44").ClearContents()
$ws.Range("H44").Value = 3500
$ws.Range("I44").Value = 3500
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 3500
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -3012

$ws.Range("H55").Value = 16053
$ws.Range("J55").Value = 16053
$ws.Range("L55").Value = 16053
$ws.Range("N55").Value = -16683

$ws.Range("H61").Value = 1970.8182
$ws.Range("I61").Value = 1967.9
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1967.9
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1755.9
$ws.Range("N61").Value = -2424

$ws.Range("H80").Value = 22553.428
$ws.Range("J80").Value = 22553.428
$ws.Range("L80").Value = 22553.428
$ws.Range("N80").Value = -24549.428

$ws.Range("H83").Value = 22553.428
$ws.Range("J83").Value = 22553.428
$ws.Range("L83").Value = 67660.284
$ws.Range("N83").Value = -77644.284

$ws.Range("H136").Value = 1970.8182
$ws.Range("I136").Value = 1967.9
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5903.700000000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3353.700000000001
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item(3)
$ws.Range("N33").ClearContents()
$ws.Range("H33").Value = 1900
$ws.Range("I33").Value = 1900
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1900
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1564

$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620

$ws.Range("H82").Value = 20247.1
$ws.Range("I82").Value = 14500
$ws.Range("J82").Value = 21683.875
$ws.Range("K82").Value = 14500
$ws.Range("L82").Value = 21683.875
$ws.Range("M82").Value = -14117
$ws.Range("N82").Value = -22449.875

$ws.Range("H85").Value = 20247.1
$ws.Range("I85").Value = 14500
$ws.Range("J85").Value = 21683.875
$ws.Range("K85").Value = 14500
$ws.Range("L85").Value = 21683.875
$ws.Range("M85").Value = -13174
$ws.Range("N85").Value = -24335.875

$ws = $wb.Worksheets.Item(4)
$ws.Range("H35").Value = 10863.75
$ws.Range("I35").Value = 1909.2858
$ws.Range("K35").Value = 1909.2858
$ws.Range("M35").Value = -1615.2858

$ws.Range("H41").Value = 12500
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -20856

$ws.Range("H50").Value = 9310.429
$ws.Range("J50").Value = 9310.429
$ws.Range("L50").Value = 9310.429
$ws.Range("N50").Value = -10560.429

$ws.Range("H51").Value = 8439.429
$ws.Range("J51").Value = 9346
$ws.Range("L51").Value = 9346
$ws.Range("N51").Value = -10818

$ws.Range("H58").Value = 683.22034
$ws.Range("I58").Value = 673.39026
$ws.Range("J58").Value = 705.6111
$ws.Range("K58").Value = 673.39026
$ws.Range("L58").Value = 705.6111
$ws.Range("M58").Value = -470.39026
$ws.Range("N58").Value = -1111.6111

$ws.Range("H60").Value = 142034
$ws.Range("J60").Value = 142034
$ws.Range("L60").Value = 142034
$ws.Range("N60").Value = -143056

$ws.Range("H61").Value = 8439.429
$ws.Range("J61").Value = 9346
$ws.Range("L61").Value = 9346
$ws.Range("N61").Value = -10042

$ws.Range("H68").Value = 18295.334
$ws.Range("J68").Value = 18295.334
$ws.Range("L68").Value = 18295.334
$ws.Range("N68").Value = -19793.334

$ws.Range("H71").Value = 18295.334
$ws.Range("J71").Value = 18295.334
$ws.Range("L71").Value = 54886.00199999999
$ws.Range("N71").Value = -62374.00199999999

$ws.Range("H109").Value = 10850
$ws.Range("J109").Value = 10971.429
$ws.Range("L109").Value = 10971.429
$ws.Range("N109").Value = -13051.429

$ws.Range("H132").Value = 2184.4187
$ws.Range("I132").Value = 1642.921
$ws.Range("J132").Value = 6299.8
$ws.Range("K132").Value = 4928.763
$ws.Range("L132").Value = 18899.4
$ws.Range("M132").Value = -2398.763
$ws.Range("N132").Value = -23959.4

$ws.Range("H134").Value = 7150553.5
$ws.Range("I134").Value = 8179.8
$ws.Range("J134").Value = 66670332
$ws.Range("K134").Value = 24539.4
$ws.Range("L134").Value = 200010996
$ws.Range("M134").Value = -22004.4
$ws.Range("N134").Value = -200016066

$ws.Range("H136").Value = 683.22034
$ws.Range("I136").Value = 673.39026
$ws.Range("J136").Value = 705.6111
$ws.Range("K136").Value = 2020.17078
$ws.Range("L136").Value = 2116.8333
$ws.Range("M136").Value = 529.8292200000001
$ws.Range("N136").Value = -7216.8333

$ws = $wb.Worksheets.Item(5)
$ws.Range("H3").Value = 10836.857
$ws.Range("I3").Value = 5353.3335
$ws.Range("J3").Value = 14949.5
$ws.Range("K3").Value = 16060.0005
$ws.Range("L3").Value = 44848.5
$ws.Range("M3").Value = -15948.0005
$ws.Range("N3").Value = -45072.5

$ws.Range("H5").Value = 834217.4
$ws.Range("I5").Value = 855.6
$ws.Range("J5").Value = 2223153.8
$ws.Range("K5").Value = 2566.8
$ws.Range("L5").Value = 6669461.399999999
$ws.Range("M5").Value = -2454.8
$ws.Range("N5").Value = -6669685.399999999

$ws.Range("H113").Value = 770.08
$ws.Range("I113").Value = 690
$ws.Range("J113").Value = 801.2222
$ws.Range("K113").Value = 2070
$ws.Range("L113").Value = 2403.6666
$ws.Range("M113").Value = 100
$ws.Range("N113").Value = -6743.6666

$ws.Range("H129").Value = 1699.125
$ws.Range("I129").Value = 902.2222
$ws.Range("J129").Value = 2177.2666
$ws.Range("K129").Value = 2706.6666
$ws.Range("L129").Value = 6531.7998
$ws.Range("M129").Value = 2293.3334
$ws.Range("N129").Value = -16531.7998

$ws.Range("H135").Value = 834217.4
$ws.Range("I135").Value = 855.6
$ws.Range("J135").Value = 2223153.8
$ws.Range("K135").Value = 7700.400000000001
$ws.Range("L135").Value = 20008384.2
$ws.Range("M135").Value = -5165.400000000001
$ws.Range("N135").Value = -20013454.2

$ws = $wb.Worksheets.Item(6)
$ws.Range("N57").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0

$ws.Range("H122").Value = 2527.318
$ws.Range("I122").Value = 2723.4443
$ws.Range("J122").Value = 2391.5386
$ws.Range("K122").Value = 8170.3329
$ws.Range("L122").Value = 7174.6158
$ws.Range("M122").Value = -5720.3329
$ws.Range("N122").Value = -12074.6158

$ws.Range("H123").Value = 34245
$ws.Range("J123").Value = 34245
$ws.Range("L123").Value = 34245
$ws.Range("N123").Value = -39145

$ws.Range("N130").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0

$ws = $wb.Worksheets.Item(7)
$ws.Range("H109").Value = 28522.143
$ws.Range("J109").Value = 28522.143
$ws.Range("L109").Value = 28522.143
$ws.Range("N109").Value = -31296.143

$ws.Range("H132").Value = 3229.4707
$ws.Range("I132").Value = 2988.6667
$ws.Range("J132").Value = 3500.375
$ws.Range("K132").Value = 8966.000100000001
$ws.Range("L132").Value = 10501.125
$ws.Range("M132").Value = -6436.000100000001
$ws.Range("N132").Value = -15561.125

$ws.Range("H136").Value = 1755.4
$ws.Range("I136").Value = 1559.3572
$ws.Range("K136").Value = 4678.071599999999
$ws.Range("M136").Value = -2128.071599999999

$ws = $wb.Worksheets.Item(8)
$ws.Range("N109").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0

$ws.Range("H127").Value = 33194
$ws.Range("J127").Value = 33194
$ws.Range("L127").Value = 33194
$ws.Range("N127").Value = -43114

$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960

$ws.Range("H132").Value = 1600.96
$ws.Range("I132").Value = 1602.6945
$ws.Range("J132").Value = 1596.5
$ws.Range("K132").Value = 4808.083500000001
$ws.Range("L132").Value = 4789.5
$ws.Range("M132").Value = -2278.083500000001
$ws.Range("N132").Value = -9849.5

$ws.Range("H136").Value = 1316.7885
$ws.Range("I136").Value = 1323
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 3969
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -1419
$ws.Range("N136").Value = -8100
